$d = $word.ActiveDocument

function Find-ParaIndex($doc, $needle) {
    $idx = 0
    $result = -1
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text.Contains($needle)) {
            $result = $idx
        }
    }
    return $result
}

# ------------------------------------------------------------------
# 1. Remove the "_GoBack" bookmark from its current location (end of
#    the "User Design" paragraph, right after "...more satisfying
#    design for this project. ").  It will be re-added later at the
#    end of the new "Implementation" paragraph.
# ------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ------------------------------------------------------------------
# 2. Locate the "Rapid Construction (Development Stage)" heading
#    paragraph and the empty paragraph right after it, then insert
#    the new "Rapid Construction" stage body paragraph, an empty
#    paragraph, and a tab-only paragraph after it (still before the
#    "Implementation (Deployment Stage)" heading).
# ------------------------------------------------------------------
$rcIdx = Find-ParaIndex $d "Rapid Construction (Development Stage)"
$afterRcIdx = $rcIdx + 1
$afterRc = $d.Paragraphs.Item($afterRcIdx)

# Insert a new (empty) paragraph right after the blank paragraph that
# follows the heading, then fill it with the stage body text.
$afterRc.Range.InsertParagraphAfter()
$devIdx = $afterRcIdx + 1
$devPara = $d.Paragraphs.Item($devIdx)
$devRange = $devPara.Range
$devRange.Collapse(1)
$devRange.InsertAfter([char]9)
$devRange.Collapse(0)
$devRange.InsertAfter("During this stage we took the web application prototype and transformed it into a working model. ")
$devRange.Collapse(0)
$devRange.InsertAfter("Although at this stage changes could still be made, due to the iterative design phase, we were able to quickly and smoothly finalize the development of this application. ")
$devRange.Collapse(0)
$devRange.InsertAfter("It was during this phase that we were able to conduct unit testing of each function, then integration testing of each functionality. Finally, we were able to conduct alpha testing of the entire we")
$devRange.Collapse(0)
$devRange.InsertAfter("b application.")

# Insert an empty paragraph, then a tab-only paragraph, after the new
# "Rapid Construction" body paragraph.
$devPara = $d.Paragraphs.Item($devIdx)
$devPara.Range.InsertParagraphAfter()
$emptyIdx = $devIdx + 1
$emptyAfterDev = $d.Paragraphs.Item($emptyIdx)
$emptyAfterDev.Range.InsertParagraphAfter()
$tabIdx = $emptyIdx + 1
$tabOnlyPara = $d.Paragraphs.Item($tabIdx)
$tabRange = $tabOnlyPara.Range
$tabRange.Collapse(1)
$tabRange.InsertAfter([char]9)

# ------------------------------------------------------------------
# 3. Locate "Implementation (Deployment Stage)" heading paragraph.
#    The blank paragraph that follows it is kept untouched; the blank
#    paragraph after *that* one is turned into the new
#    "Implementation" stage body paragraph, and the "_GoBack"
#    bookmark is re-added at its very end.
# ------------------------------------------------------------------
$implIdx = Find-ParaIndex $d "Implementation (Deployment Stage)"
$emptyAfterImplIdx = $implIdx + 1
$finalIdx = $emptyAfterImplIdx + 1
$finalPara = $d.Paragraphs.Item($finalIdx)
$finalRange = $finalPara.Range
$finalRange.Collapse(1)
$finalRange.InsertAfter([char]9)
$finalRange.Collapse(0)
$finalRange.InsertAfter("At this stage the final tests were done to confirm full functionality before finally launching the web application fully. ")
$finalRange.Collapse(0)
$finalRange.InsertAfter("Upon completion, we had a functioning application that could be used with node.js to satisfy the requirements outlined in the assignment.")

$finalRange.Collapse(0)
# `Bookmarks.Add` on a genuinely zero-length Range is unreliable, so
# insert a throw-away placeholder character, bookmark the 1-character
# range around it, then delete the placeholder -- the bookmark
# collapses to a zero-length bookmark at that position (matching the
# original "_GoBack" bookmark's shape).
$finalRange.InsertAfter("X")
$placeholderRange = $d.Range($finalRange.Start, $finalRange.Start + 1)
$d.Bookmarks.Add("_GoBack", $placeholderRange)
$placeholderRange.Text = ""
